$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds an Excel date serial. For rows 2-33 it was
# updated from 45243 (2023-11-13) to 45244 (2023-11-14).
$ws.Range("C2:C33").Value = 45244
